$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 35, shifting existing rows 35-38 down to 36-39
$ws.Rows("35:35").Insert()

# Populate the new row 35 with the latest weekly price data
$ws.Range("A35").Value = 8
$ws.Range("B35").Value = "Terminal La Palmera de La Serena"
$ws.Range("C35").Value = "Coquimbo"
$ws.Range("D35").Value = 45154
$ws.Range("D35").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E35").Value = 4
$ws.Range("F35").Value = 100112013
$ws.Range("G35").Value = "Alcachofa"
$ws.Range("H35").Value = "Española"
$ws.Range("I35").Value = "Primera"
$ws.Range("J35").Value = 500
$ws.Range("K35").Value = 11000
$ws.Range("L35").Value = 12000
$ws.Range("M35").Value = 11500
$ws.Range("N35").Value = "$/caja 30 unidades"
$ws.Range("O35").Value = "Provincia de Limarí"
$ws.Range("P35").Value = 383
$ws.Range("Q35").Value = 30
$ws.Range("R35").Value = "Hortaliza"
